$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 88654af0-...md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 18:56:43"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 88654af0-... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-13 18:56:35"
$wsZhCn.Range("K3").Value = "2016-08-13 18:57:10"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 88654af0-... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-13 18:56:43"
$wsDeDe.Range("K3").Value = "2016-08-13 18:57:20"
